$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 545, pushing the existing rows 545:631 down to 549:635.
$ws.Range("A545:R548").EntireRow.Insert()

# Row 545 (new): Larga vida / Primera, $/bandeja 18 kilos, Región de Arica y Parinacota
$ws.Range("A545").Value2 = 8
$ws.Range("B545").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C545").Value2 = "Coquimbo"
$ws.Range("D545").Value2 = 44474
$ws.Range("E545").Value2 = 4
$ws.Range("F545").Value2 = 100112020
$ws.Range("G545").Value2 = "Tomate"
$ws.Range("H545").Value2 = "Larga vida"
$ws.Range("I545").Value2 = "Primera"
$ws.Range("J545").Value2 = 800
$ws.Range("K545").Value2 = 16000
$ws.Range("L545").Value2 = 17000
$ws.Range("M545").Value2 = 16500
$ws.Range("N545").Value2 = "$/bandeja 18 kilos"
$ws.Range("O545").Value2 = "Región de Arica y Parinacota"
$ws.Range("P545").Value2 = 917
$ws.Range("Q545").Value2 = 18
$ws.Range("R545").Value2 = "Hortaliza"

# Row 546 (new): Larga vida / Primera, $/caja 10 kilos, Región de Arica y Parinacota
$ws.Range("A546").Value2 = 8
$ws.Range("B546").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C546").Value2 = "Coquimbo"
$ws.Range("D546").Value2 = 44474
$ws.Range("E546").Value2 = 4
$ws.Range("F546").Value2 = 100112020
$ws.Range("G546").Value2 = "Tomate"
$ws.Range("H546").Value2 = "Larga vida"
$ws.Range("I546").Value2 = "Primera"
$ws.Range("J546").Value2 = 760
$ws.Range("K546").Value2 = 7800
$ws.Range("L546").Value2 = 8000
$ws.Range("M546").Value2 = 7900
$ws.Range("N546").Value2 = "$/caja 10 kilos"
$ws.Range("O546").Value2 = "Región de Arica y Parinacota"
$ws.Range("P546").Value2 = 790
$ws.Range("Q546").Value2 = 10
$ws.Range("R546").Value2 = "Hortaliza"

# Row 547 (new): Larga vida / Segunda, $/bandeja 18 kilos, Región de Arica y Parinacota
$ws.Range("A547").Value2 = 8
$ws.Range("B547").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C547").Value2 = "Coquimbo"
$ws.Range("D547").Value2 = 44474
$ws.Range("E547").Value2 = 4
$ws.Range("F547").Value2 = 100112020
$ws.Range("G547").Value2 = "Tomate"
$ws.Range("H547").Value2 = "Larga vida"
$ws.Range("I547").Value2 = "Segunda"
$ws.Range("J547").Value2 = 400
$ws.Range("K547").Value2 = 14000
$ws.Range("L547").Value2 = 15000
$ws.Range("M547").Value2 = 14500
$ws.Range("N547").Value2 = "$/bandeja 18 kilos"
$ws.Range("O547").Value2 = "Región de Arica y Parinacota"
$ws.Range("P547").Value2 = 806
$ws.Range("Q547").Value2 = 18
$ws.Range("R547").Value2 = "Hortaliza"

# Row 548 (new): Larga vida / Segunda, $/caja 10 kilos, Región de Arica y Parinacota
$ws.Range("A548").Value2 = 8
$ws.Range("B548").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C548").Value2 = "Coquimbo"
$ws.Range("D548").Value2 = 44474
$ws.Range("E548").Value2 = 4
$ws.Range("F548").Value2 = 100112020
$ws.Range("G548").Value2 = "Tomate"
$ws.Range("H548").Value2 = "Larga vida"
$ws.Range("I548").Value2 = "Segunda"
$ws.Range("J548").Value2 = 480
$ws.Range("K548").Value2 = 6500
$ws.Range("L548").Value2 = 7000
$ws.Range("M548").Value2 = 6750
$ws.Range("N548").Value2 = "$/caja 10 kilos"
$ws.Range("O548").Value2 = "Región de Arica y Parinacota"
$ws.Range("P548").Value2 = 675
$ws.Range("Q548").Value2 = 10
$ws.Range("R548").Value2 = "Hortaliza"
